# "Surf/2.xlsx" update: rewording several referral-program menu strings,
# renaming "Помощь" -> "F.A.Q." and "Регистрация на поездку" -> "Забронировать",
# and adding a new trailing row for "Ваш email:".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Referral / bonus-program strings switched from "your" (formal, вашей) to
# "your" (informal, твоей) and the phone-related captions reworded.
$ws.Range("A28").Value = "Накоплено бонусов:  "
$ws.Range("A29").Value = "Людей зарегестрировалось по твоей ссылке:  "
$ws.Range("A30").Value = "Людей оплатило тур по твоей ссылке:  "
$ws.Range("A31").Value = "Сменить телефон"
$ws.Range("A32").Value = "Твой телефон:  "

# Menu labels renamed.
$ws.Range("A47").Value = "F.A.Q."
$ws.Range("A49").Value = "Забронировать"

# New trailing row (A59) with the "Ваш email:" caption, matching the
# formatting (wrap text) used by the rest of column A.
$ws.Range("A59").Value = "Ваш email:"
$ws.Range("A59").WrapText = $true

# Restore the selection/active-cell state to match the refreshed view.
$ws.Range("A32").Select() | Out-Null
